# Automatische test-sync: 2025-08-02 00:12:50
#
# Adds the newly-logged "Testmail #17" row to the Logs sheet, rolls the
# matching category tally into the Dashboard sheet, extends the
# conditional-formatting ranges + chart source ranges to cover the new
# rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 7
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A7").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Range("B7").Value = "mailmind.test@zohomail.eu"
$logs.Range("C7").Value = "Testmail #17: Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Range("D7").Value = "Planning / Afspraak"
$logs.Range("E7").Value = "Beste afzender,`nBedank u voor uw e-mail. Ik bevestig graag de afspraak voor de demo bij Van Dijk op vrijdag om 11:00 uur.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Range("F7").Value = "2025-08-02 00:12:20"
$logs.Range("G7").Value = "Ja"
$logs.Range("H7").Value = "Nee"
$logs.Range("I7").Value = "Ja"
$logs.Range("J7").Value = "Nee"

# Extend the conditional-formatting sqref for every column block from
# row 6 to row 7 (modifying one rule per block is enough, since every
# cfRule within a block shares the same sqref).
$logs.Range("D2:D6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D7"))
$logs.Range("G2:G6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G7"))
$logs.Range("H2:H6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H7"))
$logs.Range("I2:I6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I7"))
$logs.Range("J2:J6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J7"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: append row 5 (new category tally)
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A5").Value = "Planning / Afspraak"
$dash.Range("B5").Value = 1

# ---------------------------------------------------------------------
# 3. Chart: extend the category/value source ranges to include row 5
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)

$series.XValues = "='Dashboard'!`$A`$2:`$A`$5"
$series.Values = "='Dashboard'!`$B`$2:`$B`$5"
